# Updated cryptos list on Wed Dec 20 02:18:53 UTC 2023 with GitHub Actions
#
# Refresh the Coin/Price/Volume(1h) table: most rows only get new Price
# and Volume(1h) figures, but ImmutableX and InternetComputer(DFINITY)
# also swap places (rows 23/24) as their ranking changed.
#
# Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps them as text (matching the original "inline string" cells
# and preserving formatting such as trailing zeros, e.g. "58.00"),
# instead of silently converting them to numbers. ClearFormats() afterwards
# drops the "quote prefix" text-format styling Excel would otherwise stamp
# on the cell, so the cell keeps the workbook's default (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.321.73"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.179.33"
$ws.Range("E3").Value = "  -2.92%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'251.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").Value = "'0.609"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("D7").Value = "'74.37"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.27%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.576"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.69%  "
$ws.Range("D10").Value = "'39.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D13").Value = "'6.69"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.62%  "
$ws.Range("D14").Value = "2.503.54"
$ws.Range("E14").Value = "  -3.22%  "
$ws.Range("D15").Value = "'14.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.03%  "
$ws.Range("D16").Value = "2.172.58"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "'0.762"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -7.00%  "
$ws.Range("D18").Value = "42.221.20"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("E19").Value = "  -4.50%  "
$ws.Range("D20").Value = "'70.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'5.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").Value = "'225.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'2.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.91%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'9.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -15.63%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'10.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.03%  "
$ws.Range("D27").Value = "'3.40"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").Value = "'2.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.68%  "
$ws.Range("D30").Value = "'37.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'171.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "'19.95"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "'0.0815"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'5.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.80%  "
$ws.Range("D35").Value = "'0.119"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "'0.105"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.11%  "
$ws.Range("E37").Value = "  -5.74%  "
$ws.Range("D38").Value = "'0.0328"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "'11.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -12.31%  "
$ws.Range("D41").Value = "'5.13"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -8.61%  "
$ws.Range("D42").Value = "'0.192"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").Value = "'2.52"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +8.83%  "
$ws.Range("D44").Value = "'58.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.73%  "
$ws.Range("D45").Value = "'100.88"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").Value = "'0.0966"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").Value = "'8.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("E51").Value = "  -0.04%  "
